$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 picks up the values that used to live in row 3 (A2 id, D2 time)
# Force A2 to stay text (it looks numeric) by briefly using a text format,
# then restoring the cell's default style so no stray formatting sticks.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "232005"
$ws.Range("A2").Style = "Normal"

$ws.Range("D2").Value = "10:18:25"

# Row 3 (now a duplicate of the updated row 2) is removed, shifting rows up
$ws.Rows.Item(3).Delete()
